$d = $word.ActiveDocument

# The paragraph currently reads:  "ตาราง … Sequence Diagram"
# and must become:                "ตารางที่ 1 Sequence Diagram"
#
# Concretely (per the target OOXML):
#   - the run holding the single space right after "ตาราง" becomes "ที่ "
#     (and stops being rendered via the w:hint="cs" fallback, while still
#     keeping its <w:cs/> marker)
#   - the run holding the ellipsis "… " becomes "1"
#   - a new run holding a single " " is inserted right after it, so that
#     "Sequence" (and " Diagram" after it) stay in their own runs.

# --- locate "ตาราง" ---------------------------------------------------
$findTable = $d.Content
$findTable.Find.ClearFormatting()
$findTable.Find.Forward = $true
$findTable.Find.Text = "ตาราง"
$findTable.Find.Execute() | Out-Null
if (-not $findTable.Find.Found) {
    throw "edit.ps1: could not find 'ตาราง'"
}
$afterTable = $findTable.End

# --- run 1: the lone space -> "ที่ " -----------------------------------
$spaceRun = $d.Range($afterTable, $afterTable + 1)
if ($spaceRun.Text -ne " ") {
    throw "edit.ps1: expected a single space after 'ตาราง', got [$($spaceRun.Text)]"
}
$spaceRun.Text = "ที่ "

# Re-derive rFonts for the run so the w:hint="cs" fallback attribute is
# dropped (the run keeps ascii/hAnsi "TH Sarabun New" and <w:cs/>).
$retextedSpace = $d.Range($afterTable, $afterTable + 4)
$retextedSpace.Font.Name = "TH Sarabun New"

# --- run 2: the ellipsis "… " -> "1" + new " " run ---------------------
$ellipsisStart = $afterTable + 4
$ellipsisRun = $d.Range($ellipsisStart, $ellipsisStart + 2)
if ($ellipsisRun.Text -ne [char]8230 + " ") {
    throw "edit.ps1: expected ellipsis run, got [$($ellipsisRun.Text)]"
}
$ellipsisRun.Text = "1"

# Insert the separate trailing-space run right after "1".
$spaceInsertionPoint = $d.Range($ellipsisStart + 1, $ellipsisStart + 1)
$spaceInsertionPoint.InsertAfter(" ")

# Nudge formatting on the new "1" / " " runs so they do not get silently
# re-merged into the following "Sequence"/" Diagram" runs (which must stay
# as their own separate runs, matching the original document).
$rOne = $d.Range($ellipsisStart, $ellipsisStart + 1)
$rOne.Font.Bold = $true
$rOne.Font.Bold = $false

$rSpace = $d.Range($ellipsisStart + 1, $ellipsisStart + 2)
$rSpace.Font.Bold = $true
$rSpace.Font.Bold = $false

# --- keep "Sequence" / " Diagram" as distinct runs ---------------------
# (search starting right where the edited text ends, so we land on THIS
# "Sequence"/"Diagram" pair and not an earlier occurrence elsewhere in
# the document, e.g. in the title.)
$findSequence = $d.Range($ellipsisStart + 2, $d.Content.End)
$findSequence.Find.ClearFormatting()
$findSequence.Find.Forward = $true
$findSequence.Find.Text = "Sequence"
$findSequence.Find.Execute() | Out-Null
if ($findSequence.Find.Found) {
    $rSeq = $d.Range($findSequence.Start, $findSequence.End)
    $rSeq.Font.Bold = $true
    $rSeq.Font.Bold = $false

    $findDiagram = $d.Range($findSequence.End, $d.Content.End)
    $findDiagram.Find.ClearFormatting()
    $findDiagram.Find.Forward = $true
    $findDiagram.Find.Text = "Diagram"
    $findDiagram.Find.Execute() | Out-Null
    if ($findDiagram.Find.Found) {
        $rDia = $d.Range($findSequence.End, $findDiagram.End)
        $rDia.Font.Bold = $true
        $rDia.Font.Bold = $false
    }
}
